# Refresh the model's date window by one day (shift every timestamp and its
# derived "Lookup" label forward by one calendar day), matching the daily
# re-run of the Entsoe Unintended Deviation data fetch ("Horeco").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $newSerial = $aCell.Value2 + 1
    $aCell.Value2 = $newSerial

    $dCell = $ws.Cells.Item($r, 4)
    $quarter = $dCell.Value2

    $newDate = [DateTime]::FromOADate($newSerial)
    $dateText = $newDate.ToString("dd.MM.yyyy")

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value2 = "$dateText$quarter"
}
